$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: for D-column values that look numeric (single dot,
# e.g. "322.77"), Excel auto-converts them to a number on assignment. The
# source data stores these as plain text, so we force text by switching the
# cell to a text number format, assigning the value, and then restoring the
# original style/format so no stray style attribute is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
$ws.Range("D2").Value = "30.202.09"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.983.53"
$ws.Range("E3").Value = "  +6.00%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "322.77"
$ws.Range("E5").Value = "  +1.13%  "
Set-TextValue "D6" "0.9989"
$ws.Range("E6").Value = "  -0.18%  "
Set-TextValue "D7" "0.5115"
$ws.Range("E7").Value = "  +1.15%  "
Set-TextValue "D8" "0.4108"
$ws.Range("E8").Value = "  +3.64%  "
Set-TextValue "D9" "0.08439"
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("E10").Value = "  +3.83%  "
Set-TextValue "D11" "42.58"
Set-TextValue "D12" "24.18"
$ws.Range("E12").Value = "  +3.29%  "
$ws.Range("D13").Value = "1.961.35"
$ws.Range("E13").Value = "  +5.26%  "
Set-TextValue "D14" "6.472"
$ws.Range("E14").Value = "  +3.04%  "
Set-TextValue "D15" "7.396"
$ws.Range("E15").Value = "  +2.81%  "
Set-TextValue "D16" "0.9968"
$ws.Range("E16").Value = "  -0.51%  "
Set-TextValue "D17" "93.77"
$ws.Range("E17").Value = "  +2.12%  "
Set-TextValue "D18" "0.00001105"
$ws.Range("E18").Value = "  +1.85%  "
Set-TextValue "D19" "0.06544"
$ws.Range("E19").Value = "  +1.57%  "
Set-TextValue "D20" "18.78"
$ws.Range("E20").Value = "  +3.85%  "
Set-TextValue "D21" "0.9985"
$ws.Range("E21").Value = "  -0.21%  "
Set-TextValue "D22" "6.075"
$ws.Range("E22").Value = "  +3.97%  "
$ws.Range("D23").Value = "30.258.74"
$ws.Range("E23").Value = "  +0.80%  "
Set-TextValue "D24" "11.48"
$ws.Range("E24").Value = "  +3.23%  "
Set-TextValue "D25" "2.205"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "2.196.73"
$ws.Range("E26").Value = "  +5.60%  "
Set-TextValue "D27" "22.56"
$ws.Range("E27").Value = "  +6.03%  "
Set-TextValue "D28" "162.91"
$ws.Range("E28").Value = "  +1.16%  "
Set-TextValue "D29" "2.376"
$ws.Range("E29").Value = "  +7.03%  "
Set-TextValue "D30" "130.74"
$ws.Range("E30").Value = "  +2.70%  "
Set-TextValue "D31" "1.135"
$ws.Range("E31").Value = "  +5.96%  "
Set-TextValue "D32" "0.1055"
$ws.Range("E32").Value = "  +2.03%  "
Set-TextValue "D33" "6.029"
$ws.Range("E33").Value = "  +1.49%  "
Set-TextValue "D34" "3.820"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  +12.10%  "
Set-TextValue "D36" "0.02473"
$ws.Range("E36").Value = "  +1.52%  "
Set-TextValue "D37" "5.373"
$ws.Range("E37").Value = "  +2.80%  "
Set-TextValue "D38" "0.06496"
$ws.Range("E38").Value = "  +2.22%  "
Set-TextValue "D39" "0.2173"
$ws.Range("E39").Value = "  +1.51%  "
Set-TextValue "D40" "8.896"
$ws.Range("E40").Value = "  +4.91%  "
Set-TextValue "D41" "0.6587"
$ws.Range("E41").Value = "  +4.47%  "
Set-TextValue "D42" "11.81"
$ws.Range("E42").Value = "  +4.55%  "
Set-TextValue "D43" "1.221"
$ws.Range("E43").Value = "  +0.56%  "
Set-TextValue "D44" "13.54"
$ws.Range("E44").Value = "  +3.70%  "
Set-TextValue "D45" "0.6116"
$ws.Range("E45").Value = "  +3.38%  "
Set-TextValue "D46" "2.182"
$ws.Range("E46").Value = "  +4.09%  "
Set-TextValue "D47" "3.639"
$ws.Range("E47").Value = "  +0.46%  "
Set-TextValue "D50" "79.57"
$ws.Range("E50").Value = "  +2.81%  "
Set-TextValue "D51" "0.06899"
$ws.Range("E51").Value = "  +2.06%  "

# --- Row 48/49: EOS and Quant swap positions, with updated price/volume ---
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D48" "123.82"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D49" "1.221"
$ws.Range("E49").Value = "  +1.32%  "
